# Updates cryptos price/volume columns (D, E) for rows 2-51 per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper scratch cell used to force plain-numeric-looking strings (e.g. "588.17")
# to be written as TEXT instead of being auto-coerced to a number by Excel,
# while leaving the destination cell's style/number-format untouched
# (format-as-text, copy, paste-special-values, then clear the scratch cell).
$scratch = $ws.Range("Z1")

$ws.Range("D2").Value = '60.791.06'
$ws.Range("E2").Value = '  -3.56%  '
$ws.Range("D3").Value = '2.908.59'
$ws.Range("E3").Value = '  -4.06%  '
$ws.Range("E4").Value = '  +0.01%  '
$scratch.NumberFormat = "@"
$scratch.Value = '588.17'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E5").Value = '  -0.99%  '
$scratch.NumberFormat = "@"
$scratch.Value = '144.30'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E6").Value = '  -6.37%  '
$ws.Range("E7").Value = '  +0.03%  '
$scratch.NumberFormat = "@"
$scratch.Value = '0.503'
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E8").Value = '  -2.40%  '
$ws.Range("D9").Value = '2.908.60'
$scratch.NumberFormat = "@"
$scratch.Value = '6.72'
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E10").Value = '  -3.53%  '
$ws.Range("E11").Value = '  -5.48%  '
$scratch.NumberFormat = "@"
$scratch.Value = '0.444'
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E12").Value = '  -4.32%  '
$ws.Range("E13").Value = '  -3.94%  '
$scratch.NumberFormat = "@"
$scratch.Value = '33.44'
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E14").Value = '  -6.20%  '
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").Value = '3.387.55'
$ws.Range("E16").Value = '  -4.13%  '
$ws.Range("D17").Value = '60.747.58'
$ws.Range("E17").Value = '  -3.56%  '
$scratch.NumberFormat = "@"
$scratch.Value = '6.68'
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E18").Value = '  -5.86%  '
$ws.Range("D19").Value = '2.908.26'
$ws.Range("E19").Value = '  -4.05%  '
$scratch.NumberFormat = "@"
$scratch.Value = '427.83'
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E20").Value = '  -5.78%  '
$scratch.NumberFormat = "@"
$scratch.Value = '13.52'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E21").Value = '  -5.49%  '
$scratch.NumberFormat = "@"
$scratch.Value = '0.682'
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E22").Value = '  -2.39%  '
$scratch.NumberFormat = "@"
$scratch.Value = '7.07'
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E23").Value = '  -6.15%  '
$scratch.NumberFormat = "@"
$scratch.Value = '80.78'
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E24").Value = '  -2.93%  '
$scratch.NumberFormat = "@"
$scratch.Value = '10.80'
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E25").Value = '  -6.32%  '
$ws.Range("E26").Value = '  -5.76%  '
$scratch.NumberFormat = "@"
$scratch.Value = '11.88'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E27").Value = '  -4.49%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  +0.01%  '
$scratch.NumberFormat = "@"
$scratch.Value = '2.20'
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E30").Value = '  -3.42%  '
$scratch.NumberFormat = "@"
$scratch.Value = '7.20'
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E31").Value = '  -4.19%  '
$ws.Range("E32").Value = '  -3.78%  '
$scratch.NumberFormat = "@"
$scratch.Value = '26.42'
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E33").Value = '  -4.55%  '
$ws.Range("E34").Value = '  -3.71%  '
$ws.Range("D35").Value = '0.0₃0857'
$ws.Range("E35").Value = '  -1.03%  '
$ws.Range("E36").Value = '  -2.64%  '
$scratch.NumberFormat = "@"
$scratch.Value = '5.58'
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E37").Value = '  -5.69%  '
$ws.Range("E38").Value = '  -4.61%  '
$scratch.NumberFormat = "@"
$scratch.Value = '49.40'
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("E40").Value = '  -5.82%  '
$ws.Range("E41").Value = '  -5.96%  '
$ws.Range("E42").Value = '  -5.75%  '
$ws.Range("E43").Value = '  -4.99%  '
$scratch.NumberFormat = "@"
$scratch.Value = '41.29'
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E44").Value = '  -7.60%  '
$scratch.NumberFormat = "@"
$scratch.Value = '0.0350'
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E45").Value = '  -3.03%  '
$scratch.NumberFormat = "@"
$scratch.Value = '372.19'
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E46").Value = '  -5.79%  '
$ws.Range("D47").Value = '2.695.38'
$ws.Range("E47").Value = '  -0.96%  '
$scratch.NumberFormat = "@"
$scratch.Value = '132.36'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E48").Value = '  -0.40%  '
$scratch.NumberFormat = "@"
$scratch.Value = '24.08'
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("E50").Value = '  -6.57%  '
$ws.Range("E51").Value = '  -2.97%  '
